$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column A: Experiment names (rows 2-18) ----
$ws.Range("A2").Value  = "Off1"
$ws.Range("A3").Value  = "Off2"
$ws.Range("A4").Value  = "Off3"
$ws.Range("A5").Value  = "Off4"
$ws.Range("A6").Value  = "Off5"
$ws.Range("A7").Value  = "Off6"
$ws.Range("A8").Value  = "Bygd1"
$ws.Range("A9").Value  = "Bygd2"
$ws.Range("A10").Value = "Bygd3"
$ws.Range("A11").Value = "Bygd4"
$ws.Range("A12").Value = "Bygd5"
$ws.Range("A13").Value = "Bygd6"
$ws.Range("A14").Value = "Bygd7"
$ws.Range("A15").Value = "Fjell1"
$ws.Range("A16").Value = "Fjell2"
$ws.Range("A17").Value = "Fjell3"
$ws.Range("A18").Value = "Fjell4"

# ---- Column B: span (rows 5-18 newly filled in) ----
$ws.Range("B5").Value  = 135
$ws.Range("B6").Value  = 180
$ws.Range("B7").Value  = 225
$ws.Range("B8").Value  = 15
$ws.Range("B9").Value  = 45
$ws.Range("B10").Value = 90
$ws.Range("B11").Value = 135
$ws.Range("B12").Value = 180
$ws.Range("B13").Value = 225
$ws.Range("B14").Value = 270
$ws.Range("B15").Value = 90
$ws.Range("B16").Value = 135
$ws.Range("B17").Value = 180
$ws.Range("B18").Value = 225

# ---- Column C: distMiil (rows 8-18 updated) ----
$ws.Range("C8").Value  = 8
$ws.Range("C9").Value  = 8
$ws.Range("C10").Value = 8
$ws.Range("C11").Value = 8
$ws.Range("C12").Value = 8
$ws.Range("C13").Value = 8
$ws.Range("C14").Value = 8
$ws.Range("C15").Value = 6
$ws.Range("C16").Value = 6
$ws.Range("C17").Value = 6
$ws.Range("C18").Value = 6

# ---- Column D: Overlap (various rows) ----
$ws.Range("D5").Value  = "1, 2, 4, 5, 6"
$ws.Range("D6").Value  = "1, 2, 3, 4, 5, 6"
$ws.Range("D7").Value  = "1, 2, 3, 4, 5, 6, 7, 8, 9"
# D8 becomes a blank cell that merely carries the quote-prefix style
# (matches the author typing a leading apostrophe and leaving the cell
# otherwise empty). We first set it so Excel applies the quotePrefix
# style, then clear the value back out again.
$ws.Range("D8").Value  = "'"
$ws.Range("D8").Value  = ""
# D10 holds the text "2" (not the number 2). A leading apostrophe forces
# Excel to store it as text and apply the quote-prefix style, exactly
# like D2/D3 above.
$ws.Range("D10").Value = "'2"
$ws.Range("D11").Value = "2, 6"
$ws.Range("D12").Value = "2, 3, 6"
$ws.Range("D13").Value = "2, 3, 6, 8"
$ws.Range("D14").Value = "2, 3, 6, 8"

# ---- Column E: # overlap (rows 5-18 newly filled in) ----
$ws.Range("E5").Value  = 5
$ws.Range("E6").Value  = 6
$ws.Range("E7").Value  = 9
$ws.Range("E8").Value  = 0
$ws.Range("E9").Value  = 0
$ws.Range("E10").Value = 1
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 3
$ws.Range("E13").Value = 4
$ws.Range("E14").Value = 4
$ws.Range("E15").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("E18").Value = 0

# Leave column F (Success rate) alone -- it is a shared formula
# (=E/G*100) and recalculates automatically from the new E/G values.

# Move the active selection to E19, matching where the author left off.
$excel.Goto($ws.Range("E19"))
